$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test case data rows (26 and 27), following the existing pattern and
# matching the original authoring order so new shared strings are appended in
# the same sequence: TB2485530, 11678849, TN2485531, 11678853.

# Row 26
$ws.Range("A26").Value = "TB2485530"

# The Quote Number values look numeric; write them via TEXT()+paste-values so
# they land as plain text cells (matching the rest of the column) without
# leaving behind any new cell style.
$ws.Range("B26").Formula = '=TEXT(11678849,"0")'
$ws.Range("B26").Copy()
$ws.Range("B26").PasteSpecial(-4163)

$ws.Range("C26").Value = "Bond - No Credit"
$ws.Range("D26").Value = "Nil"
$ws.Range("E26").Value = "TC005"

# Row 27
$ws.Range("A27").Value = "TN2485531"

$ws.Range("B27").Formula = '=TEXT(11678853,"0")'
$ws.Range("B27").Copy()
$ws.Range("B27").PasteSpecial(-4163)

$ws.Range("C27").Value = "Personal Auto - Credit"
$ws.Range("D27").Value = "Base"
$ws.Range("E27").Value = "TC001"

$excel.CutCopyMode = 0
